$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1891025641025641
$ws.Range("C2").Value = 0.5801282051282052
$ws.Range("J2").Value = 0.009615384615384616
$ws.Range("P2").Value = 0.1474358974358974
$ws.Range("S2").Value = 0.07371794871794872
# Row 3
$ws.Range("B3").Value = 0.02604166666666667
$ws.Range("C3").Value = 0.02604166666666667
$ws.Range("J3").Value = 0.02604166666666667
$ws.Range("P3").Value = 0.734375
$ws.Range("S3").Value = 0.1875
# Row 4
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.5945945945945946
$ws.Range("S4").Value = 0.3783783783783784
# Row 6
$ws.Range("B6").Value = 0.06557377049180328
$ws.Range("D6").Value = 0.004098360655737705
$ws.Range("F6").Value = 0.06967213114754098
$ws.Range("J6").Value = 0.2786885245901639
$ws.Range("O6").Value = 0.05327868852459016
$ws.Range("Q6").Value = 0.1434426229508197
$ws.Range("R6").Value = 0.05327868852459016
$ws.Range("S6").Value = 0.3319672131147541
# Row 7
$ws.Range("B7").Value = 0.1371681415929203
$ws.Range("D7").Value = 0.01769911504424779
$ws.Range("F7").Value = 0.06637168141592921
$ws.Range("J7").Value = 0.1327433628318584
$ws.Range("O7").Value = 0.008849557522123894
$ws.Range("Q7").Value = 0.1238938053097345
$ws.Range("R7").Value = 0.07079646017699115
$ws.Range("S7").Value = 0.4424778761061947
# Row 8
$ws.Range("B8").Value = 0.1092233009708738
$ws.Range("D8").Value = 0.01941747572815534
$ws.Range("F8").Value = 0.08009708737864078
$ws.Range("J8").Value = 0.1359223300970874
$ws.Range("O8").Value = 0.02912621359223301
$ws.Range("Q8").Value = 0.1359223300970874
$ws.Range("R8").Value = 0.09951456310679611
$ws.Range("S8").Value = 0.3907766990291262
# Row 9
$ws.Range("B9").Value = 0.07446808510638298
$ws.Range("D9").Value = 0.01595744680851064
$ws.Range("F9").Value = 0.101063829787234
$ws.Range("J9").Value = 0.148936170212766
$ws.Range("O9").Value = 0.03191489361702127
$ws.Range("Q9").Value = 0.1223404255319149
$ws.Range("R9").Value = 0.09042553191489362
$ws.Range("S9").Value = 0.4148936170212766
# Row 10
$ws.Range("B10").Value = 0.1159305993690852
$ws.Range("D10").Value = 0.02050473186119874
$ws.Range("E10").Value = 0.0007886435331230284
$ws.Range("F10").Value = 0.06545741324921135
$ws.Range("J10").Value = 0.1317034700315458
$ws.Range("O10").Value = 0.0220820189274448
$ws.Range("Q10").Value = 0.1782334384858044
$ws.Range("R10").Value = 0.09621451104100946
$ws.Range("S10").Value = 0.3690851735015773
# Row 11
$ws.Range("G11").Value = 0.1549707602339181
$ws.Range("J11").Value = 0.09064327485380116
$ws.Range("K11").Value = 0.195906432748538
$ws.Range("L11").Value = 0.5497076023391813
$ws.Range("S11").Value = 0.008771929824561403
# Row 12
$ws.Range("G12").Value = 0.7658536585365854
$ws.Range("J12").Value = 0.1707317073170732
$ws.Range("K12").Value = 0.00975609756097561
$ws.Range("L12").Value = 0.01463414634146342
$ws.Range("S12").Value = 0.03902439024390244
# Row 13
$ws.Range("G13").Value = 0.6341463414634146
$ws.Range("J13").Value = 0.2439024390243902
$ws.Range("S13").Value = 0.1219512195121951
# Row 15
$ws.Range("F15").Value = 0.03149606299212598
$ws.Range("H15").Value = 0.1102362204724409
$ws.Range("I15").Value = 0.07874015748031496
$ws.Range("J15").Value = 0.3543307086614173
$ws.Range("K15").Value = 0.07086614173228346
$ws.Range("M15").Value = 0.003937007874015748
$ws.Range("O15").Value = 0.07086614173228346
$ws.Range("S15").Value = 0.2795275590551181
# Row 16
$ws.Range("F16").Value = 0.02985074626865672
$ws.Range("H16").Value = 0.1990049751243781
$ws.Range("I16").Value = 0.04975124378109453
$ws.Range("J16").Value = 0.3781094527363184
$ws.Range("K16").Value = 0.1691542288557214
$ws.Range("M16").Value = 0.02487562189054726
$ws.Range("O16").Value = 0.02487562189054726
$ws.Range("S16").Value = 0.1243781094527363
# Row 17
$ws.Range("F17").Value = 0.0215633423180593
$ws.Range("H17").Value = 0.1617250673854447
$ws.Range("I17").Value = 0.08355795148247978
$ws.Range("J17").Value = 0.4070080862533693
$ws.Range("K17").Value = 0.1078167115902965
$ws.Range("M17").Value = 0.0215633423180593
$ws.Range("O17").Value = 0.07008086253369272
$ws.Range("S17").Value = 0.1266846361185984
# Row 18
$ws.Range("F18").Value = 0.01456310679611651
$ws.Range("H18").Value = 0.1747572815533981
$ws.Range("I18").Value = 0.1019417475728155
$ws.Range("J18").Value = 0.4223300970873786
$ws.Range("K18").Value = 0.08737864077669903
$ws.Range("M18").Value = 0.02427184466019417
$ws.Range("N18").Value = 0.004854368932038835
$ws.Range("O18").Value = 0.05339805825242718
$ws.Range("S18").Value = 0.116504854368932
# Row 19
$ws.Range("F19").Value = 0.02641802641802642
$ws.Range("H19").Value = 0.1872571872571873
$ws.Range("I19").Value = 0.07303807303807304
$ws.Range("J19").Value = 0.358974358974359
$ws.Range("K19").Value = 0.1227661227661228
$ws.Range("M19").Value = 0.01631701631701632
$ws.Range("N19").Value = 0.002331002331002331
$ws.Range("O19").Value = 0.08080808080808081
$ws.Range("S19").Value = 0.1320901320901321

Write-Host "Applied 109 cell updates"
